# Generate Report for Handback
# The handback run completed for both e2e test files: the file that was
# previously "Ready for handoff" (a8acdbae-212f-4577-a97a-bcb26a9fcc8c) has
# now been handed back in sync with en-US, same as c2de2976-cfd8-49c4-966d-c086ca65102a.
# The status report table is regenerated with the two rows swapped (a8acdbae
# first, c2de2976 second) and fresh handback timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md"
$ovw.Range("B2").Value = "Handed back: in sync with en-US"
$ovw.Range("C2").Value = "Handed back: in sync with en-US"

$ovw.Range("A3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.md"
$ovw.Range("B3").Value = "Handed back: in sync with en-US"
$ovw.Range("C3").Value = "Handed back: in sync with en-US"

$ovw.Range("A4").Value = ".localization-config"
$ovw.Range("B4").Value = "Not to be localized"
$ovw.Range("C4").Value = "Not to be localized"

$ovw.Hyperlinks.Delete()
$ovw.Hyperlinks.Add($ovw.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md")
$ovw.Hyperlinks.Add($ovw.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.md")
$ovw.Hyperlinks.Add($ovw.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf"
$zh.Range("D2").Value = "2016-02-22 14:32:35"
$zh.Range("E2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md"
$zh.Range("F2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf"
$zh.Range("G2").Value = "2016-02-22 14:33:25"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf"
$zh.Range("D3").Value = "2016-02-22 14:32:35"
$zh.Range("E3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.md"
$zh.Range("F3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf"
$zh.Range("G3").Value = "2016-02-22 14:33:25"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/138554b0e6326bc41202dd941ec0c49f655f0bd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/mt/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e7568d8c9dbac9b492ad59a8c9895812c4ff3f19/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/11f9a4c5ea6ab85a0223af967a740320c12cd8cf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/138554b0e6326bc41202dd941ec0c49f655f0bd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/mt/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e7568d8c9dbac9b492ad59a8c9895812c4ff3f19/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/11f9a4c5ea6ab85a0223af967a740320c12cd8cf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf"
$de.Range("D2").Value = "2016-02-22 14:32:48"
$de.Range("E2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md"
$de.Range("F2").Value = "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf"
$de.Range("G2").Value = "2016-02-22 14:33:48"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf"
$de.Range("D3").Value = "2016-02-22 14:32:48"
$de.Range("E3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.md"
$de.Range("F3").Value = "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf"
$de.Range("G3").Value = "2016-02-22 14:33:48"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ac32645c8b63ff09feb696e53fae8feb5cb95f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/mt/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7acca1f69aaea064609e9a8bc886803547e81867/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6c6e04908eb66992d7dbab48b9b74486b94d1bbd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf", "", "", "a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ac32645c8b63ff09feb696e53fae8feb5cb95f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/mt/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7acca1f69aaea064609e9a8bc886803547e81867/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6c6e04908eb66992d7dbab48b9b74486b94d1bbd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf", "", "", "c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2a64f690e803e9cb023fea4a951d32d8d0629f88/.localization-config", "", "", ".localization-config")
